# Added NB_SIGNALS=4 for HA_VWAP.
#
# The "TestCases" sheet had 6 sample rows (Test# 1-6): BTCUSDT/ETHUSDT each
# run with NB_SIGNALS 1/2/3, plus a handful of still-blank template rows
# below them. This edit trims the sheet back down to a single kept sample
# row (Test #1, NB_SIGNALS=1, unchanged) and a brand-new second row that
# exercises NB_SIGNALS=4 (with DistVWAP_PCT reset to 0.0), then drops the
# now-redundant rows so the blank template rows collapse back down to a
# single row beneath the data, matching the sheet's original "one spare
# row" shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# --- trim the trailing blank template rows down to just one (old 8:14 -> new 4:8) ---
$ws.Rows("13:14").Delete() | Out-Null

# --- drop the old Test #2/#3 (BTCUSDT) and #4/#5/#6 (ETHUSDT) sample rows;
#     row 8's blank template shifts up to become the new row 4 onward ---
$ws.Rows("4:7").Delete() | Out-Null

# --- renumber the kept sample row and give it the new NB_SIGNALS=4 settings ---
$ws.Range("A3").Value() = 1
$ws.Range("K3").Value() = '{"EMA": 200, "DistVWAP_PCT": 0.0, "NB_SIGNALS": 4}'

# --- restore the view state (scroll position + active cell) ---
$ws.Range("I13").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
